$d = $word.ActiveDocument

# --- Step 1: Insert a page break at the very end of the document, keeping the
# break run inside the SAME paragraph as the preceding text run. InsertBreak
# normally spawns a brand new paragraph for the break, so we immediately
# delete the paragraph mark it introduces to merge the break back into the
# original paragraph. ---
$end = $d.Content.End
$r = $d.Range($end - 1, $end - 1)
$r.InsertBreak(7)  # wdPageBreak = 7

$markPos = $end - 1
$markRange = $d.Range($markPos, $markPos + 1)
$markRange.Delete()

# --- Step 2: Append the trailing empty run (matches the source XML's final
# <w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r>) to that very same paragraph.
# InsertParagraphAfter gives us a new paragraph whose sole run carries that
# default rPr/no-text shape, then we merge the resulting paragraph mark away
# so the run ends up inside the page-break paragraph instead of starting a
# new one. ---
$end2 = $d.Content.End
$tail2 = $d.Range($end2 - 1, $end2 - 1)
$tail2.InsertParagraphAfter()

$markPos2 = $end2 - 1
$markRange2 = $d.Range($markPos2, $markPos2 + 1)
$markRange2.Delete()

# --- Step 3: Append the new stand-up entries for 09/08/22, one paragraph per
# line (blank strings become the blank spacer paragraphs seen in the diff). ---
$lines = @(
    "09/08/22",
    "Devin",
    "Yesterday Worked on page content and styling and started working on backend. Today working on the backend. Blockers - lack of knowledge",
    "",
    "Hodan",
    "Yesterday worked on page content and styling, discussion board. Today working on payment system. Blockers - backend not set up yet.",
    "",
    "Toseef",
    "Yesterday worked on page content and styling, made the email us function on the contact page. Working on backend. Blockers - lack of knowledge",
    "",
    "Waseem",
    "Yesterday on page content and styling. Today working on documentation. Blockers - "
)

foreach ($line in $lines) {
    $curEnd = $d.Content.End
    $tail = $d.Range($curEnd - 1, $curEnd - 1)
    $tail.InsertParagraphAfter()

    if ($line -ne "") {
        $newEnd = $d.Content.End
        $newPara = $d.Range($newEnd - 1, $newEnd - 1)
        $newPara.InsertAfter($line)
    }
}
